# Update crypto price/volume data as scraped on Wed Apr 19 10:13:10 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.383.86"
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -2.25%  "

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.996.10"
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -5.73%  "

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "330.86"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -4.44%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.04%  "

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4932"
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -5.11%  "

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4196"
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -5.84%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "53.03"
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -1.80%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.08859"
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -5.38%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.115"
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -5.60%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.179.97"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +3.78%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "Solana"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "23.29"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -7.94%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.055"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -5.74%  "

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.503"
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -6.72%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -6.51%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -0.11%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -4.88%  "

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06633"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -0.91%  "

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.78"
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -8.07%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +0.12%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.969"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -4.97%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.427.66"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -2.22%  "

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.84"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -6.93%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.282"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -1.44%  "

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.263.83"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -3.84%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.606"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +0.11%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "157.35"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -3.18%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.54"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -7.00%  "

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.350"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -7.14%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "127.34"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -4.94%  "

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.057"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -8.12%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.09911"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -6.18%  "

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.568"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -12.03%  "

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.849"
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -6.49%  "

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.772"
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -4.92%  "

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.580"
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -10.77%  "

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.02448"
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -6.46%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06364"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -7.38%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.283"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -3.32%  "

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.76"
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -7.35%  "

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6505"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -8.44%  "

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2071"
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -7.62%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +0.05%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6316"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -7.83%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "NEARProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.206"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -7.13%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.35"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -9.09%  "

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.262"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -0.34%  "

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.540"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -2.44%  "

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00000000338"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -3.41%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "ThetaToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.153"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -2.83%  "
